$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status of bug #10 (row 11) from OPEN to FIXED
$ws.Range("C11").Value = "FIXED"

# Add new bug #11 in row 12
$ws.Range("A12").Value = 11
$ws.Range("D12").Value = "Jobs"
$ws.Range("B12").Value = "Deleting a job with invoices throws an exception"
$ws.Range("C12").Value = "OPEN"
$ws.Range("E12").Value = "If a job has any invoices associated with it and the user attempts to delete it an expeption is thrown and reported. A more friendly error message would be better."
$ws.Range("F12").Value = 40266
$ws.Range("G12").Value = 40266

# Copy formatting from row 11 to row 12 (style, row height)
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)

$ws.Rows.Item(12).RowHeight = 45

$ws.Range("C11").Select()
